$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-11 down to 8-12
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new "PCA(6)" result
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Lasso Regression+normalization+ lag1 +PCA(6)"
$ws.Cells.Item(7, 3).Value = 88.856083384300604

# Renumber the Id column for the rows that shifted down
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10

# Copy formatting from the row above (row 6) into the new row 7 so that
# borders/styles match the rest of the table (the bottom, thick-bordered
# row keeps its own formatting automatically since it shifted with the
# insert).
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Widen column B to fit the longer model name (stored width should end up
# as exactly 46 characters)
$ws.Columns.Item(2).ColumnWidth = 45.285714285714285

# Update the active selection
$ws.Range("B14").Select()
